$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.729.65"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'3.487.94"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'592.32"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'171.61"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "'4.090.38"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "'28.82"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "'66.747.06"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "'3.490.46"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "'392.78"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "'7.91"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'72.72"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "'10.16"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'23.69"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'7.32"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "'163.15"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "'1.89"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "'6.88"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").Value = "'27.13"
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").Value = "'26.11"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "'2.793.07"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'42.66"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("D47").Value = "'335.39"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("D48").Value = "'34.28"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("E50").Value = "  -1.83%  "
